# Add data for 2022-06-04: roll the "through" date forward by one day
# (2022-05-26 -> 2022-05-27) and bump the June / Total figures that changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet tab to reflect the new "through" date.
$ws.Name = "Through 2022-05-27"

# Update the header label in I1 ("2022 (through 05-26)" -> "2022 (through 05-27)").
$ws.Range("I1").Value = "2022 (through 05-27)"

# June total (row 6) gained one incident.
$ws.Range("I6").Value = 96

# Grand total (row 14) reflects the same increment.
$ws.Range("I14").Value = 647
